$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 621
$ws.Range("J38").Value = 3500
$ws.Range("L38").Value = 10500
$ws.Range("N38").Value = -11244
$ws.Range("H39").Value = 190.03334
$ws.Range("I39").Value = 196.44827
$ws.Range("K39").Value = 589.3448100000001
$ws.Range("M39").Value = -293.3448100000001
$ws.Range("H43").Value = 12000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12138
$ws.Range("H97").Value = 1510.2142
$ws.Range("J97").Value = 1510.2142
$ws.Range("L97").Value = 4530.642599999999
$ws.Range("N97").Value = -5522.642599999999
$ws.Range("H107").Value = 532.6667
$ws.Range("I107").Value = 583.26666
$ws.Range("K107").Value = 583.26666
$ws.Range("M107").Value = 1336.73334
$ws.Range("H112").Value = 2135.6296
$ws.Range("J112").Value = 2135.6296
$ws.Range("L112").Value = 6406.888800000001
$ws.Range("N112").Value = -8622.888800000001
$ws.Range("H116").Value = 8026.875
$ws.Range("I116").Value = 7435.5
$ws.Range("J116").Value = 8618.25
$ws.Range("K116").Value = 7435.5
$ws.Range("L116").Value = 8618.25
$ws.Range("M116").Value = -3993.5
$ws.Range("N116").Value = -15502.25
$ws.Range("H138").Value = 4715.1055
$ws.Range("I138").Value = 2424.7273
$ws.Range("J138").Value = 5262.804
$ws.Range("K138").Value = 7274.1819
$ws.Range("L138").Value = 15788.412
$ws.Range("M138").Value = -2134.1819
$ws.Range("N138").Value = -26068.412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10891.171
$ws.Range("I32").Value = 9013.842000000001
$ws.Range("K32").Value = 9013.842000000001
$ws.Range("M32").Value = -8726.842000000001
$ws.Range("H61").Value = 8494.294
$ws.Range("I61").Value = 8872
$ws.Range("K61").Value = 8872
$ws.Range("M61").Value = -8660
$ws.Range("H63").Value = 5832.909
$ws.Range("I63").Value = 2360.3333
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 2360.3333
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -1674.3333
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 5832.909
$ws.Range("I66").Value = 2360.3333
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 11801.6665
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -8369.666499999999
$ws.Range("N66").Value = -56864
$ws.Range("H74").Value = 4553
$ws.Range("I74").Value = 3569
$ws.Range("J74").Value = 5537
$ws.Range("K74").Value = 3569
$ws.Range("L74").Value = 5537
$ws.Range("M74").Value = -2695
$ws.Range("N74").Value = -7285
$ws.Range("H77").Value = 4553
$ws.Range("I77").Value = 3569
$ws.Range("J77").Value = 5537
$ws.Range("K77").Value = 17845
$ws.Range("L77").Value = 27685
$ws.Range("M77").Value = -13477
$ws.Range("N77").Value = -36421
$ws.Range("H86").Value = 39997.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 39997.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H110").Value = 4216.1665
$ws.Range("I110").Value = 3063.7222
$ws.Range("J110").Value = 7673.5
$ws.Range("K110").Value = 3063.7222
$ws.Range("L110").Value = 7673.5
$ws.Range("M110").Value = -1018.7222
$ws.Range("N110").Value = -11763.5
$ws.Range("H132").Value = 4259.3
$ws.Range("I132").Value = 3721.7693
$ws.Range("K132").Value = 11165.3079
$ws.Range("M132").Value = -8635.3079
$ws.Range("H136").Value = 8494.294
$ws.Range("I136").Value = 8872
$ws.Range("K136").Value = 26616
$ws.Range("M136").Value = -24066

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 39999.5
$ws.Range("J60").Value = 39999.5
$ws.Range("L60").Value = 39999.5
$ws.Range("N60").Value = -41021.5
$ws.Range("H132").Value = 3193.7678
$ws.Range("I132").Value = 2761.9807
$ws.Range("K132").Value = 8285.9421
$ws.Range("M132").Value = -5755.9421

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 199.625
$ws.Range("I33").Value = 251.4
$ws.Range("J33").Value = 113.333336
$ws.Range("K33").Value = 1508.4
$ws.Range("L33").Value = 680.000016
$ws.Range("M33").Value = -1225.4
$ws.Range("N33").Value = -1246.000016
$ws.Range("H107").Value = 2232969
$ws.Range("I107").Value = 556.5
$ws.Range("J107").Value = 7814000
$ws.Range("K107").Value = 1669.5
$ws.Range("L107").Value = 23442000
$ws.Range("M107").Value = 250.5
$ws.Range("N107").Value = -23445840
$ws.Range("H122").Value = 2582.0334
$ws.Range("I122").Value = 878.7
$ws.Range("J122").Value = 3433.7
$ws.Range("K122").Value = 7908.3
$ws.Range("L122").Value = 30903.3
$ws.Range("M122").Value = -5458.3
$ws.Range("N122").Value = -35803.3
$ws.Range("H137").Value = 253572.5
$ws.Range("J137").Value = 253572.5
$ws.Range("L137").Value = 760717.5
$ws.Range("N137").Value = -770917.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2568.8
$ws.Range("J6").Value = 2933.5
$ws.Range("L6").Value = 2933.5
$ws.Range("N6").Value = -3159.5
$ws.Range("H13").Value = 580.8182
$ws.Range("I13").Value = 841.6667
$ws.Range("J13").Value = 267.8
$ws.Range("K13").Value = 841.6667
$ws.Range("L13").Value = 267.8
$ws.Range("M13").Value = -702.6667
$ws.Range("N13").Value = -545.8
$ws.Range("H16").Value = 2568.8
$ws.Range("J16").Value = 2933.5
$ws.Range("L16").Value = 2933.5
$ws.Range("N16").Value = -3433.5
$ws.Range("H24").Value = 30061.555
$ws.Range("J24").Value = 36388.75
$ws.Range("L24").Value = 36388.75
$ws.Range("N24").Value = -36734.75
$ws.Range("H27").Value = 4000
$ws.Range("I27").Value = 3000
$ws.Range("K27").Value = 3000
$ws.Range("M27").Value = -2834
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H54").Value = 10666.667
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -11780
$ws.Range("H126").Value = 5043.143
$ws.Range("I126").Value = 6081.4
$ws.Range("J126").Value = 4264.45
$ws.Range("K126").Value = 18244.2
$ws.Range("L126").Value = 12793.35
$ws.Range("M126").Value = -15774.2
$ws.Range("N126").Value = -17733.35
$ws.Range("H132").Value = 4098.625
$ws.Range("I132").Value = 1962.5
$ws.Range("K132").Value = 5887.5
$ws.Range("M132").Value = -3357.5
$ws.Range("H134").Value = 57146.875
$ws.Range("J134").Value = 57146.875
$ws.Range("L134").Value = 171440.625
$ws.Range("N134").Value = -176510.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8362.208000000001
$ws.Range("I7").Value = 6014.65
$ws.Range("K7").Value = 6014.65
$ws.Range("M7").Value = -5902.65
$ws.Range("H46").Value = 4327.4546
$ws.Range("I46").Value = 1533.3334
$ws.Range("J46").Value = 5375.25
$ws.Range("K46").Value = 1533.3334
$ws.Range("L46").Value = 5375.25
$ws.Range("M46").Value = -1345.3334
$ws.Range("N46").Value = -5751.25
$ws.Range("H61").Value = 14501
$ws.Range("I61").Value = 12000.6
$ws.Range("K61").Value = 12000.6
$ws.Range("M61").Value = -11798.6
$ws.Range("H69").Value = 38499.5
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
$ws.Range("H72").Value = 38499.5
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112
$ws.Range("H113").Value = 14501
$ws.Range("I113").Value = 12000.6
$ws.Range("K113").Value = 12000.6
$ws.Range("M113").Value = -9830.6
$ws.Range("H126").Value = 8362.208000000001
$ws.Range("I126").Value = 6014.65
$ws.Range("K126").Value = 18043.95
$ws.Range("M126").Value = -15573.95
$ws.Range("H132").Value = 5364.2964
$ws.Range("I132").Value = 4633.24
$ws.Range("K132").Value = 13899.72
$ws.Range("M132").Value = -11369.72
$ws.Range("H136").Value = 5867.9165
$ws.Range("I136").Value = 4102
$ws.Range("K136").Value = 12306
$ws.Range("M136").Value = -9756

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -26040
$ws.Range("H113").Value = 288.37036
$ws.Range("J113").Value = 145.75
$ws.Range("L113").Value = 437.25
$ws.Range("N113").Value = -4777.25
$ws.Range("H132").Value = 6151
$ws.Range("I132").Value = 3832
$ws.Range("J132").Value = 12335
$ws.Range("K132").Value = 11496
$ws.Range("L132").Value = 37005
$ws.Range("M132").Value = -8966
$ws.Range("N132").Value = -42065
